# Update "want to go" counts (column F) on the 展览 (sheet 1) and
# 全部类型 (sheet 4) worksheets, which list the same events but at
# slightly different row offsets (sheet 4 contains two extra rows).
#
# Sheet "展览" (index 1): F4, F11, F13, F27, F37, F42
# Sheet "全部类型" (index 4): F4, F11, F13, F28, F38, F44

$wb = $excel.ActiveWorkbook

$wsExhibit = $wb.Worksheets.Item(1)   # 展览
$wsAll     = $wb.Worksheets.Item(4)   # 全部类型

# 展览 sheet updates
$wsExhibit.Range("F4").Value  = 7655
$wsExhibit.Range("F11").Value = 4427
$wsExhibit.Range("F13").Value = 118
$wsExhibit.Range("F27").Value = 1451
$wsExhibit.Range("F37").Value = 3118
$wsExhibit.Range("F42").Value = 820

# 全部类型 sheet updates (same events, shifted rows)
$wsAll.Range("F4").Value  = 7655
$wsAll.Range("F11").Value = 4427
$wsAll.Range("F13").Value = 118
$wsAll.Range("F28").Value = 1451
$wsAll.Range("F38").Value = 3118
$wsAll.Range("F44").Value = 820
